# Updates cryptocurrency price (D) and 1h volume change (E) columns
# per the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.686.03"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").Value = "1.960.71"
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("E5").Value = "  +0.90%  "

$ws.Range("E6").Value = "  +1.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.61"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.375"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.58%  "

$ws.Range("E10").Value = "  -6.57%  "

$ws.Range("E11").Value = "  +0.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.834"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.92%  "

$ws.Range("D15").Value = "2.245.36"
$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("E16").Value = "  +3.12%  "

$ws.Range("D17").Value = "1.954.32"
$ws.Range("E17").Value = "  +0.59%  "

$ws.Range("D18").Value = "36.567.41"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.88%  "

$ws.Range("E20").Value = "  -1.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "230.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.05%  "

$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("E24").Value = "  +4.75%  "

$ws.Range("E25").Value = "  +2.52%  "

$ws.Range("E26").Value = "  +7.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.74%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.04%  "

$ws.Range("E30").Value = "  +17.79%  "

$ws.Range("E31").Value = "  +1.45%  "

$ws.Range("E32").Value = "  +4.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0615"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("E34").Value = "  +7.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.92%  "

$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0980"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.56%  "

$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("E42").Value = "  +2.27%  "

$ws.Range("E43").Value = "  +0.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.01%  "

$ws.Range("D45").Value = "1.369.97"
$ws.Range("E45").Value = "  +2.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.51%  "

$ws.Range("E47").Value = "  +1.12%  "

$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.62%  "

$ws.Range("D51").Value = "2.136.80"
$ws.Range("E51").Value = "  +0.95%  "
